# Apply "adding averages and more checks" update to the Training Dashboard sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Training Dashboard"

# --- 1. Update PERIOD TO EXPIRE (col H) and LAST UPDATE (col I) for rows 3-27 ---
# Every row's "last update" date moves from 08-Sep-2025 to 16-Sep-2025 (8 days later),
# which reduces the "period to expire" count by 8 for every row.

$newPeriod = @{
    3  = 470
    4  = 490
    5  = 447
    6  = 531
    7  = 483
    8  = 483
    9  = 219
    10 = 503
    11 = 541
    12 = 491
    13 = 210
    14 = 335
    15 = 335
    16 = 706
    17 = 360
    18 = 323
    19 = -103
    20 = -180
    21 = -41
    22 = -41
    23 = 176
    24 = 155
    25 = 268
    26 = 323
    27 = 348
}

# Make sure column I keeps being a literal text date (not auto-converted to a
# real Excel date serial number) by formatting it as Text before writing.
$ws.Range("I3:I27").NumberFormat = "@"

foreach ($row in 3..27) {
    $ws.Cells.Item($row, 8).Value = $newPeriod[$row]     # column H
    $ws.Cells.Item($row, 9).Value = "16-Sep-2025"        # column I
}

# --- 2. Header / title font color -> bold white text ---
# The title (A1) and the table header row (A2:K2) now share the same bold,
# white font (the header keeps its dark-blue fill and borders). These two
# cell styles are shared across the whole workbook, so the "Exam Dashboard"
# sheet's title/header are updated the same way.
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215      # RGB(255,255,255) = white
$ws.Range("A2:K2").Font.Color = 16777215   # RGB(255,255,255) = white

$ws2 = $wb.Worksheets.Item(2)   # "Exam Dashboard"
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215
